$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit the new, longer Khmer/English clip-name text.
$ws.Columns.Item(3).ColumnWidth = 39

# Add a new use-case row (row 16): "to confirm a good input" -> "OK"
$ws.Range("A16").Value = "to confirm a good input"
$ws.Range("B16").Value = "OK"
$ws.Range("C16").Value = "OK"
$ws.Range("I16").Value = "still to add, suggested by bong Thavy"
$ws.Rows.Item(16).RowHeight = 30

# Rows 5 and 6 no longer need the taller wrapped height.
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15

# Update the view: zoom to 100% and move the active selection.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("I17").Select()
